$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # ALC
$ws2 = $wb.Worksheets.Item(2)  # ARM
$ws3 = $wb.Worksheets.Item(3)  # BSM
$ws4 = $wb.Worksheets.Item(4)  # CRP
$ws5 = $wb.Worksheets.Item(5)  # CUL
$ws6 = $wb.Worksheets.Item(6)  # GSM
$ws7 = $wb.Worksheets.Item(7)  # LTW
$ws8 = $wb.Worksheets.Item(8)  # WVR

# Row from hunk @ 2130 (ALC)
$ws1.Range("H31").Value = 1993
$ws1.Range("I31").Value = 1993
$ws1.Range("K31").Value = 5979
$ws1.Range("M31").Value = -5749

# Row from hunk @ 6409 (ALC)
$ws1.Range("H116").Value = 25800536
$ws1.Range("J116").Value = 6065.3335
$ws1.Range("L116").Value = 6065.3335
$ws1.Range("N116").Value = -12949.3335

# Row from hunk @ 6651 (ALC)
$ws1.Range("H121").Value = 4663.278
$ws1.Range("J121").Value = 4663.278
$ws1.Range("L121").Value = 13989.834
$ws1.Range("N121").Value = -17483.834

# Row from hunk @ 7337 (ALC)
$ws1.Range("H135").Value = 2350.5483
$ws1.Range("I135").Value = 670.86957
$ws1.Range("K135").Value = 6037.826129999999
$ws1.Range("M135").Value = -3502.826129999999

# Row from hunk @ 7438 (ALC)
$ws1.Range("H137").Value = 11448849
$ws1.Range("I137").Value = 835256.7
$ws1.Range("J137").Value = 18524578
$ws1.Range("K137").Value = 2505770.1
$ws1.Range("L137").Value = 55573734
$ws1.Range("M137").Value = -2503220.1
$ws1.Range("N137").Value = -55578834

# Row from hunk @ 7643 (ALC)
$ws1.Range("H141").Value = 2205.5715
$ws1.Range("I141").Value = 2205.5715
$ws1.Range("K141").Value = 6616.7145
$ws1.Range("M141").Value = -1436.7145

# Row from hunk @ 7786 (ARM)
$ws2.Range("H2").Value = 649010.5
$ws2.Range("I2").Value = 796090.2
$ws2.Range("K2").Value = 796090.2
$ws2.Range("M2").Value = -795977.2

# Row from hunk @ 9214 (ARM)
$ws2.Range("H32").Value = 19626.438
$ws2.Range("I32").Value = 22059.486
$ws2.Range("K32").Value = 22059.486
$ws2.Range("M32").Value = -21772.486

# Row from hunk @ 9361 (ARM)
$ws2.Range("H35").Value = 3800
$ws2.Range("I35").Value = 3800
$ws2.Range("K35").Value = 3800
$ws2.Range("M35").Value = -3394

# Row from hunk @ 9845 (ARM)
$ws2.Range("H45").Value = 3090
$ws2.Range("I45").Value = 1832.5
$ws2.Range("J45").Value = 4599
$ws2.Range("K45").Value = 1832.5
$ws2.Range("L45").Value = 4599
$ws2.Range("M45").Value = -1455.5
$ws2.Range("N45").Value = -5353

# Row from hunk @ 10611 (ARM)
$ws2.Range("H61").Value = 8920.120000000001
$ws2.Range("I61").Value = 9095.380999999999
$ws2.Range("K61").Value = 9095.380999999999
$ws2.Range("M61").Value = -8883.380999999999

# Row from hunk @ 11239 (ARM)
$ws2.Range("H74").Value = 1336.2916
$ws2.Range("J74").Value = 1998.4286
$ws2.Range("L74").Value = 1998.4286
$ws2.Range("N74").Value = -3746.4286

# Row from hunk @ 11383 (ARM)
$ws2.Range("H77").Value = 1336.2916
$ws2.Range("J77").Value = 1998.4286
$ws2.Range("L77").Value = 9992.143
$ws2.Range("N77").Value = -18728.143

# Row from hunk @ 13270 (ARM)
$ws2.Range("H116").Value = 649010.5
$ws2.Range("I116").Value = 796090.2
$ws2.Range("K116").Value = 796090.2
$ws2.Range("M116").Value = -793796.2

# Row from hunk @ 13371 (ARM)
$ws2.Range("H118").Value = 72000
$ws2.Range("I118").Value = 44000
$ws2.Range("K118").Value = 44000
$ws2.Range("M118").Value = -42343

# Row from hunk @ 14045 (ARM)
$ws2.Range("H132").Value = 13569.095
$ws2.Range("I132").Value = 17714.656
$ws2.Range("K132").Value = 53143.96799999999
$ws2.Range("M132").Value = -50613.96799999999

# Row from hunk @ 14192 (ARM)
$ws2.Range("H135").Value = 74482.25
$ws2.Range("J135").Value = 74482.25
$ws2.Range("L135").Value = 74482.25
$ws2.Range("N135").Value = -84622.25

# Row from hunk @ 14241 (ARM)
$ws2.Range("H136").Value = 8920.120000000001
$ws2.Range("I136").Value = 9095.380999999999
$ws2.Range("K136").Value = 27286.143
$ws2.Range("M136").Value = -24736.143

# Row from hunk @ 14684 (BSM)
$ws3.Range("H3").Value = 649010.5
$ws3.Range("I3").Value = 796090.2
$ws3.Range("K3").Value = 796090.2
$ws3.Range("M3").Value = -795976.2

# Row from hunk @ 15600 (BSM)
$ws3.Range("H22").Value = 428.53333
$ws3.Range("I22").Value = 531.125
$ws3.Range("K22").Value = 531.125
$ws3.Range("M22").Value = -358.125

# Row from hunk @ 20998 (BSM)
$ws3.Range("H134").Value = 1366
$ws3.Range("I134").Value = 1342.0625
$ws3.Range("J134").Value = 1557.5
$ws3.Range("K134").Value = 4026.1875
$ws3.Range("L134").Value = 4672.5
$ws3.Range("M134").Value = -1491.1875
$ws3.Range("N134").Value = -9742.5

# Row from hunk @ 21524 (CRP)
$ws4.Range("H3").Value = 0
$ws4.Range("J3").Value = 0
$ws4.Range("L3").Value = 0
$ws4.Range("N3").ClearContents()

# Row from hunk @ 22262 (CRP)
$ws4.Range("H18").Value = 0
$ws4.Range("I18").Value = 0
$ws4.Range("K18").Value = 0
$ws4.Range("M18").ClearContents()

# Row from hunk @ 22363 (CRP)
$ws4.Range("H20").Value = 690000
$ws4.Range("J20").Value = 690000
$ws4.Range("L20").Value = 690000
$ws4.Range("N20").Value = -690472

# Row from hunk @ 22749 (CRP)
$ws4.Range("H28").Value = 40000
$ws4.Range("J28").Value = 40000
$ws4.Range("L28").Value = 40000
$ws4.Range("N28").Value = -40490

# Row from hunk @ 22844 (CRP)
$ws4.Range("H30").Value = 690000
$ws4.Range("J30").Value = 690000
$ws4.Range("L30").Value = 690000
$ws4.Range("N30").Value = -690182

# Row from hunk @ 22893 (CRP)
$ws4.Range("H31").Value = 5123.75
$ws4.Range("I31").Value = 2075.9
$ws4.Range("J31").Value = 5649.241
$ws4.Range("K31").Value = 2075.9
$ws4.Range("L31").Value = 5649.241
$ws4.Range("M31").Value = -1780.9
$ws4.Range("N31").Value = -6239.241

# Row from hunk @ 23043 (CRP)
$ws4.Range("H34").Value = 5123.75
$ws4.Range("I34").Value = 2075.9
$ws4.Range("J34").Value = 5649.241
$ws4.Range("K34").Value = 2075.9
$ws4.Range("L34").Value = 5649.241
$ws4.Range("M34").Value = -1873.9
$ws4.Range("N34").Value = -6053.241

# Row from hunk @ 24213 (CRP)
$ws4.Range("H58").Value = 354332.97
$ws4.Range("I58").Value = 418438
$ws4.Range("J58").Value = 307711.12
$ws4.Range("K58").Value = 418438
$ws4.Range("L58").Value = 307711.12
$ws4.Range("M58").Value = -418235
$ws4.Range("N58").Value = -308117.12

# Row from hunk @ 26596 (CRP)
$ws4.Range("H107").Value = 1818819.2
$ws4.Range("I107").Value = 2273299
$ws4.Range("K107").Value = 2273299
$ws4.Range("M107").Value = -2271379

# Row from hunk @ 27622 (CRP)
$ws4.Range("H128").Value = 690000
$ws4.Range("J128").Value = 690000
$ws4.Range("L128").Value = 690000
$ws4.Range("N128").Value = -699960

# Row from hunk @ 27812 (CRP)
$ws4.Range("H132").Value = 8340708.5
$ws4.Range("I132").Value = 9532000
$ws4.Range("K132").Value = 28596000
$ws4.Range("M132").Value = -28593470

# Row from hunk @ 27913 (CRP)
$ws4.Range("H134").Value = 2585.0908
$ws4.Range("I134").Value = 2778.6
$ws4.Range("J134").Value = 650
$ws4.Range("K134").Value = 8335.799999999999
$ws4.Range("L134").Value = 1950
$ws4.Range("M134").Value = -5800.799999999999
$ws4.Range("N134").Value = -7020

# Row from hunk @ 28014 (CRP)
$ws4.Range("H136").Value = 354332.97
$ws4.Range("I136").Value = 418438
$ws4.Range("J136").Value = 307711.12
$ws4.Range("K136").Value = 1255314
$ws4.Range("L136").Value = 923133.36
$ws4.Range("M136").Value = -1252764
$ws4.Range("N136").Value = -928233.36

# Row from hunk @ 31548 (CUL)
$ws5.Range("H64").Value = 12030.8
$ws5.Range("J64").Value = 14788.5
$ws5.Range("L64").Value = 44365.5
$ws5.Range("N64").Value = -44905.5

# Row from hunk @ 31701 (CUL)
$ws5.Range("H67").Value = 12030.8
$ws5.Range("J67").Value = 14788.5
$ws5.Range("L67").Value = 44365.5
$ws5.Range("N67").Value = -46237.5

# Row from hunk @ 33672 (CUL)
$ws5.Range("H106").Value = 23749.25
$ws5.Range("J106").Value = 23749.25
$ws5.Range("L106").Value = 71247.75
$ws5.Range("N106").Value = -73139.75

# Row from hunk @ 34131 (CUL)
$ws5.Range("H115").Value = 6918.2856
$ws5.Range("I115").Value = 2732
$ws5.Range("J115").Value = 12500
$ws5.Range("K115").Value = 8196
$ws5.Range("L115").Value = 37500
$ws5.Range("M115").Value = -7021
$ws5.Range("N115").Value = -39850

# Row from hunk @ 41843 (GSM)
$ws6.Range("H132").Value = 605688.3
$ws6.Range("I132").Value = 289882.44
$ws6.Range("K132").Value = 869647.3200000001
$ws6.Range("M132").Value = -867117.3200000001

# Row from hunk @ 43101 (LTW)
$ws7.Range("H16").Value = 2882.3447
$ws7.Range("I16").Value = 1609.091
$ws7.Range("J16").Value = 6884
$ws7.Range("K16").Value = 1609.091
$ws7.Range("L16").Value = 6884
$ws7.Range("M16").Value = -1439.091
$ws7.Range("N16").Value = -7224

# Row from hunk @ 44547 (LTW)
$ws7.Range("H46").Value = 5226.727
$ws7.Range("I46").Value = 2375
$ws7.Range("J46").Value = 5620.069
$ws7.Range("K46").Value = 2375
$ws7.Range("L46").Value = 5620.069
$ws7.Range("M46").Value = -2187
$ws7.Range("N46").Value = -5996.069

# Row from hunk @ 47142 (LTW)
$ws7.Range("H100").Value = 12044.4
$ws7.Range("I100").Value = 6000
$ws7.Range("K100").Value = 6000
$ws7.Range("M100").Value = -5459

# Row from hunk @ 52132 (WVR)
$ws8.Range("H62").Value = 7488.615
$ws8.Range("I62").Value = 5158.1665
$ws8.Range("J62").Value = 9486.143
$ws8.Range("K62").Value = 5158.1665
$ws8.Range("L62").Value = 9486.143
$ws8.Range("M62").Value = -4534.1665
$ws8.Range("N62").Value = -10734.143

# Row from hunk @ 52279 (WVR)
$ws8.Range("H65").Value = 7488.615
$ws8.Range("I65").Value = 5158.1665
$ws8.Range("J65").Value = 9486.143
$ws8.Range("K65").Value = 25790.8325
$ws8.Range("L65").Value = 47430.715
$ws8.Range("M65").Value = -22670.8325
$ws8.Range("N65").Value = -53670.715

# Row from hunk @ 53054 (WVR)
$ws8.Range("H81").Value = 2319167
$ws8.Range("I81").Value = 2976929.8
$ws8.Range("K81").Value = 5953859.6
$ws8.Range("M81").Value = -5952798.6

# Row from hunk @ 53204 (WVR)
$ws8.Range("H84").Value = 2319167
$ws8.Range("I84").Value = 2976929.8
$ws8.Range("K84").Value = 29769298
$ws8.Range("M84").Value = -29763994

# Row from hunk @ 55015 (WVR)
$ws8.Range("H122").Value = 3623.652
$ws8.Range("I122").Value = 3623.652
$ws8.Range("K122").Value = 10870.956
$ws8.Range("M122").Value = -8420.956

# Row from hunk @ 55499 (WVR)
$ws8.Range("H132").Value = 22731578
$ws8.Range("J132").Value = 71439430
$ws8.Range("L132").Value = 214318290
$ws8.Range("M132").Value = -214323350
